# SyntheticDataPipeline: use both CmsCollection and SchemaCollection
#
# The pipeline used to write collection "hasPart"/work data only into the
# SchemaCollection sheet. It now also emits a CmsCollection sheet (the
# schema:Collection-shaped rows move there, keyed by @id/hasPart/image/title),
# while SchemaCollection is trimmed back down to just its "@graph" header.

$wb = $excel.ActiveWorkbook

# --- Add the new "CmsCollection" worksheet, placed right after
#     "SchemaCollection" (the last sheet in the workbook). ---
$schemaCollection = $wb.Worksheets.Item("SchemaCollection")
$cmsCollection = $wb.Worksheets.Add($null, $schemaCollection)
$cmsCollection.Name = "CmsCollection"

# Header row.
$cmsCollection.Range("A1").Value = "@id"
$cmsCollection.Range("B1").Value = "hasPart"
$cmsCollection.Range("C1").Value = "image"
$cmsCollection.Range("D1").Value = "title"

# Data row (moved over from SchemaCollection, with an updated image
# thumbnail reference and hasPart left blank).
$cmsCollection.Range("A2").Value = "http://example.com/collection1"
$cmsCollection.Range("B2").Value = ""
$cmsCollection.Range("C2").Value = "http://example.com/collection1:Image1:Thumbnail400x400"
$cmsCollection.Range("D2").Value = "Collection1"

# --- Trim "SchemaCollection" back down to just its "@graph" header cell. ---
$schemaCollection.Range("A2:D2").ClearContents()
